# Tracking Progress update - Aug 7th
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Track Table")

# A3: simplify block name
$ws.Range("A3").Value = "MOSBius Interface "

# Row 4 (Ring Oscillator): schematic status moves to PVT/MC in process
$ws.Range("I4").Value = "PVT/MC in process"

# Row 5 (XOR Phase Detector): now in progress, start date moved up, schematic status PVT/MC in process
$ws.Range("D5").Value = "In progress"
$ws.Range("E5").Value = 45867
$ws.Range("I5").Value = "PVT/MC in process"

# Row 6 (Comparator): now in progress, start date moved up, schematic status TT in process
$ws.Range("D6").Value = "In progress"
$ws.Range("E6").Value = 45867
$ws.Range("I6").Value = "TT in process"

# Row 7 (Edge detector): now in progress, start date moved up, schematic status PVT/MC in process
$ws.Range("D7").Value = "In progress"
$ws.Range("E7").Value = 45874
$ws.Range("I7").Value = "PVT/MC in process"

# Row 8 (Integrator): now in progress, start date moved up, schematic status PVT/MC in process
$ws.Range("D8").Value = "In progress"
$ws.Range("E8").Value = 45869
$ws.Range("I8").Value = "PVT/MC in process"

# Row 9 (Low Pass Filter): now in progress, start date pushed back, schematic status PVT/MC in process
$ws.Range("D9").Value = "In progress"
$ws.Range("E9").Value = 45910
$ws.Range("I9").Value = "PVT/MC in process"

# Row 10 (Neuron): now in progress, start date pushed back, schematic status TT in process
$ws.Range("D10").Value = "In progress"
$ws.Range("E10").Value = 45910
$ws.Range("I10").Value = "TT in process"
